# Atualizações referentes a documentação e login
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update "Classificação" (column C) for rows 14 and 16 from "Essencial" to "Importante"
$ws.Range("C14").Value = "Importante"
$ws.Range("C16").Value = "Importante"

# Update the sheet's view: scroll back to top (remove topLeftCell override) and
# change the active selection to H1:I6
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$ws.Range("H1:I6").Select() | Out-Null

Write-Output "done"
